$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Medellin"
$ws.Range("B3").Value = "Cali"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1

$ws.Range("A4").Value = "Bogota"
$ws.Range("B4").Value = "Cucuta"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0

$ws.Range("F4").Select()
